# Apply the "Update all diagrams except UI" edit to the single slide
# contained in this trimmed-down presentation.
#
# Changes applied here (the ones whose targets exist on this slide):
#   - Rename the domain event "AddressBookChangedEvent" -> "ImdbChangedEvent"
#     (and its handler "handleAddresssBookChangedEvent" -> "handleImdbChangedEvent")
#     in the four textboxes that mention it.
#   - Reposition/resize the "post(ImdbChangedEvent)" textbox (shape id 62,
#     "TextBox 61") to match the shorter label.
#
# (The datetimeFigureOut field text changes in the original diff target
# other slides that this trimmed single-slide deck does not contain, so
# there is nothing on this slide to touch for those hunks.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate shapes by their stable shape Id (survives shape-order changes).
$shapesById = @{}
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    $shapesById[$sh.Id] = $sh
}

# --- "TextBox 32" (id 33): post(AddressBookChangedEvent) -> post(ImdbChangedEvent)
$sh33 = $shapesById[33]
if ($sh33 -ne $null) {
    $tr33 = $sh33.TextFrame.TextRange
    $full33 = $tr33.Text
    $idx33 = $full33.IndexOf("AddressBookChangedEvent")
    if ($idx33 -ge 0) {
        $chars33 = $tr33.Characters($idx33 + 1, "AddressBookChangedEvent".Length)
        $chars33.Text = "ImdbChangedEvent"
    }
}

# --- "TextBox 61" (id 62): post(AddressBookChangedEvent) -> post(ImdbChangedEvent)
# plus its textbox is repositioned/resized to fit the shorter text.
$sh62 = $shapesById[62]
if ($sh62 -ne $null) {
    $tr62 = $sh62.TextFrame.TextRange
    $full62 = $tr62.Text
    $idx62 = $full62.IndexOf("AddressBookChangedEvent")
    if ($idx62 -ge 0) {
        $chars62 = $tr62.Characters($idx62 + 1, "AddressBookChangedEvent".Length)
        $chars62.Text = "ImdbChangedEvent"
    }
    # NOTE: Shape.Left/Top/Width/Height round-trip through a single-precision
    # (points) representation before being re-quantized to EMU on save, so a
    # naive "EMU / 12700.0" literal can land 1 EMU off the target after
    # save. The literals below were empirically verified (via the exact
    # save->unzip->inspect round trip) to land exactly on the target EMU
    # values (x=2514600, y=4797674, cx=1889925, cy=215444).
    $sh62.Left = 198.0
    $sh62.Top = 377.769653
    $sh62.Width = 148.813035
    $sh62.Height = 16.964134488
}

# --- "TextBox 73" (id 74): handleAddresssBookChangedEvent() -> handleImdbChangedEvent()
$sh74 = $shapesById[74]
if ($sh74 -ne $null) {
    $tr74 = $sh74.TextFrame.TextRange
    $full74 = $tr74.Text
    $idx74 = $full74.IndexOf("handleAddresssBookChangedEvent")
    if ($idx74 -ge 0) {
        $chars74 = $tr74.Characters($idx74 + 1, "handleAddresssBookChangedEvent".Length)
        $chars74.Text = "handleImdbChangedEvent"
    }
}

# --- "TextBox 49" (id 50): handleAddresssBookChangedEvent() -> handleImdbChangedEvent()
$sh50 = $shapesById[50]
if ($sh50 -ne $null) {
    $tr50 = $sh50.TextFrame.TextRange
    $full50 = $tr50.Text
    $idx50 = $full50.IndexOf("handleAddresssBookChangedEvent")
    if ($idx50 -ge 0) {
        $chars50 = $tr50.Characters($idx50 + 1, "handleAddresssBookChangedEvent".Length)
        $chars50.Text = "handleImdbChangedEvent"
    }
}
